$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update sheet name (Through 2022-03-10 -> Through 2022-03-11)
$ws.Name = "Through 2022-03-11"

# Update column header text for March 2022 column
$ws.Range("B1").Value = "March 2022 (through March 11)"

# Neighborhood label shift for rows 9-12 (Lake View / Chicago Lawn / West Loop / Englewood)
$ws.Range("A9").Value = "Chicago Lawn"
$ws.Range("A10").Value = "West Loop"
$ws.Range("A11").Value = "Englewood"
$ws.Range("A12").Value = "Lake View"

# Update numeric cell values
$ws.Range("B3").Value = 3
$ws.Range("E4").Value = 3
$ws.Range("B5").Value = 3
$ws.Range("B9").Value = 3
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 1
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 2
$ws.Range("M9").Value = 1
$ws.Range("O9").Value = 2
$ws.Range("P9").Value = 4
$ws.Range("R9").Value = 2
$ws.Range("U9").Value = 4
$ws.Range("W9").Value = 1
$ws.Range("C10").Value = 4
$ws.Range("D10").Value = 2
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 2
$ws.Range("P10").Value = 2
$ws.Range("S10").Value = 1
$ws.Range("V10").Value = 2
$ws.Range("B11").Value = 6
$ws.Range("D11").Value = 12
$ws.Range("E11").Value = 1
$ws.Range("G11").Value = 4
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = 5
$ws.Range("L11").Value = 2
$ws.Range("M11").Value = 5
$ws.Range("N11").Value = 2
$ws.Range("O11").Value = 1
$ws.Range("R11").Value = 3
$ws.Range("T11").Value = 2
$ws.Range("U11").Value = 2
$ws.Range("V11").Value = 5
$ws.Range("W11").Value = 1
$ws.Range("Y11").Value = 3
$ws.Range("D12").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 2
$ws.Range("V12").Value = 1
$ws.Range("Y12").Value = 1
$ws.Range("Q17").Value = 1
$ws.Range("T17").Value = 2
$ws.Range("Q18").Value = 1
$ws.Range("W22").Value = 1
$ws.Range("N26").Value = 2
$ws.Range("K27").Value = 1
$ws.Range("N30").Value = 1
$ws.Range("N36").Value = 1
$ws.Range("Q76").Value = 2
$ws.Range("E84").Value = 2

# Clear cells that no longer have data
$ws.Range("F9").ClearContents()
$ws.Range("V9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("B10").ClearContents()
$ws.Range("E10").ClearContents()
$ws.Range("L10").ClearContents()
$ws.Range("U10").ClearContents()
$ws.Range("W10").ClearContents()
$ws.Range("F11").ClearContents()
$ws.Range("S11").ClearContents()
$ws.Range("B12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("J12").ClearContents()
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("O12").ClearContents()
$ws.Range("P12").ClearContents()
$ws.Range("R12").ClearContents()
$ws.Range("T12").ClearContents()
$ws.Range("U12").ClearContents()
$ws.Range("W12").ClearContents()
